# Actualización automática 2025-10-31 16:30:08
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" — per-product-group sales figures for row 4
# (BELTRAN ESPINOZA SONIA SARITA), row 10 (CULMA OVIEDO NINI JOHANA),
# row 21 (TAMAYO VILLACIS EDWIN XAVIER) and row 22 (TOSCANO RAMIREZ MONICA
# CECILIA), plus the "N de 24" non-zero-count summary row 26.
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("G4").Value = 122.22
$wsGrupo.Range("H4").Value = 213.3
$wsGrupo.Range("M4").Value = 2114.46
$wsGrupo.Range("N4").Value = 75.53
$wsGrupo.Range("P4").Value = 23.37

$wsGrupo.Range("L10").Value = 82.28

$wsGrupo.Range("M21").Value = 291.92

$wsGrupo.Range("M22").Value = 11501.83

$wsGrupo.Range("G26").Value = "1 de 24"
$wsGrupo.Range("H26").Value = "2 de 24"
$wsGrupo.Range("L26").Value = "3 de 24"
$wsGrupo.Range("N26").Value = "1 de 24"
$wsGrupo.Range("P26").Value = "3 de 24"

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL" — "octubre" column (F) totals per client, plus the
# column-total row 26.
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F4").Value = 2548.88
$wsMensual.Range("F10").Value = 4220.84
$wsMensual.Range("F21").Value = 842.55
$wsMensual.Range("F22").Value = 13113.85
$wsMensual.Range("F26").Value = 39199.39

# ---------------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" — VENTA / POR CUMPLIR / CUMPLIMIENTO columns
# for the groups affected by the above updates, plus the TOTAL row 14.
# Column E (5) widens from 23 to 24 characters; ColumnWidth as set via COM
# is offset by +5/6 relative to the stored OOXML character width, so we
# back that out to land exactly on 24.
# ---------------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumpl.Columns.Item(5).ColumnWidth = 23.166666666666668

# GRIFERIAS
$wsCumpl.Range("D5").Value = 122.22
$wsCumpl.Range("E5").Value = 27.78
$wsCumpl.Range("F5").Value = 0.8148

# INODOROS
$wsCumpl.Range("D6").Value = 855.91
$wsCumpl.Range("E6").Value = 51.25610861560108
$wsCumpl.Range("F6").Value = 0.9434986513177598

# NO RESURTIBLES
$wsCumpl.Range("D8").Value = 1124.63
$wsCumpl.Range("E8").Value = -458.0431724318521
$wsCumpl.Range("F8").Value = 1.687147050449365

# PIEDRA SINTERIZADA
$wsCumpl.Range("D11").Value = 4857.23
$wsCumpl.Range("E11").Value = -1935.00541814726
$wsCumpl.Range("F11").Value = 1.662168619812388

# PORCELANATO
$wsCumpl.Range("D12").Value = 27045.61
$wsCumpl.Range("E12").Value = 909.369999999999
$wsCumpl.Range("F12").Value = 0.9674701967234461

# PUERTAS DE SEGURIDAD
$wsCumpl.Range("D13").Value = 75.53
$wsCumpl.Range("E13").Value = 72.528220160454
$wsCumpl.Range("F13").Value = 0.510137160355882

# TOTAL
$wsCumpl.Range("D14").Value = 38006.5
$wsCumpl.Range("E14").Value = 4196.881100094683
$wsCumpl.Range("F14").Value = 0.9005558087836409
